# Apply the diff to the bike_store_report workbook
$wb = $excel.ActiveWorkbook

# --- Sheet "Order Status": add a new "status_text" column (D) ---
$wsOrder = $wb.Worksheets.Item("Order Status")

$wsOrder.Range("D1").Value = "status_text"
# Copy C1's formatting (bold white font on blue fill, centered, bordered)
# onto D1 so it reuses the same header cell style.
$wsOrder.Range("C1").Copy()
$wsOrder.Range("D1").PasteSpecial(-4122) | Out-Null

$wsOrder.Range("D2").Value = "Completed"
$wsOrder.Range("D3").Value = "Processing"
$wsOrder.Range("D4").Value = "Pending"
$wsOrder.Range("D5").Value = "Rejected"

# 13 units of width once the engine's +0.8333 padding is re-added on save.
$wsOrder.Columns.Item(4).ColumnWidth = 12.166666666666666

# Re-apply the autofilter so its ref grows to include the new column.
$wsOrder.Range("A1:C5").AutoFilter() | Out-Null
$wsOrder.Range("A1:D5").AutoFilter(1) | Out-Null

# The hidden _FilterDatabase defined name doesn't follow AutoFilter
# automatically here, so update it explicitly.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Order Status!_FilterDatabase") {
        $n.RefersTo = "='Order Status'!`$A`$1:`$D`$5"
    }
}

# --- Sheet "Top Products": swap rows 7 and 8 (Trek Slash <-> Electra Girl's Hawaii 20-inch) ---
$wsTop = $wb.Worksheets.Item("Top Products")

$wsTop.Range("A7").Value = "Electra Girl's Hawaii 1 (20-inch) - 2015/2016"
$wsTop.Range("B7").Value = "Electra"

$wsTop.Range("A8").Value = "Trek Slash 8 27.5 - 2016"
$wsTop.Range("B8").Value = "Trek"

# --- Sheet "Brand Prices": reorder rows 2-9 ---
$wsBrand = $wb.Worksheets.Item("Brand Prices")

$wsBrand.Range("A2").Value = "Sun Bicycles"
$wsBrand.Range("B2").Value = 529.4085567010309

$wsBrand.Range("A3").Value = "Ritchey"
$wsBrand.Range("B3").Value = 749.99

$wsBrand.Range("A4").Value = "Electra"
$wsBrand.Range("B4").Value = 519.3942278773858

$wsBrand.Range("A5").Value = "Haro"
$wsBrand.Range("B5").Value = 629.2627272727273

$wsBrand.Range("A6").Value = "Trek"
$wsBrand.Range("B6").Value = 2766.872591093118

$wsBrand.Range("A7").Value = "Surly"
$wsBrand.Range("B7").Value = 1165.934198347107

$wsBrand.Range("A8").Value = "Pure Cycles"
$wsBrand.Range("B8").Value = 441.9457364341085

$wsBrand.Range("A9").Value = "Heller"
$wsBrand.Range("B9").Value = 1400.042164948454
